# "Gimbal Lock ppp final"
#
# Slide 12 ("Wiederholung: Gimbal Lock"):
#  - Title: drop the leading "Wiederholung: " so it just reads "Gimbal Lock".
#  - Content placeholder: replace the "Gimbal Lock" line with three short
#    bullet-style paragraphs (separated by blank paragraphs).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# --- Title shape -----------------------------------------------------
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$prefix = "Wiederholung: "
$titleRange.Characters(1, $prefix.Length).Text = ""

# --- Content placeholder ----------------------------------------------
$content = $s.Shapes.Item(2)
$contentRange = $content.TextFrame.TextRange
$oldWord = "Gimbal"
# Clear just the "Gimbal" run first so the replacement text inherits the
# plain (non spell-flagged) run formatting of the remaining " Lock" run
# instead of the err="1" formatting that was attached to "Gimbal".
$contentRange.Characters(1, $oldWord.Length).Text = ""
$contentRange.Text = "drei konzentrische Ringe`r`rRing = Achse`r`rVerlust Freiheitsgrad"
